$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 1).Value = "ni Liberal Party President at Senador Mar Roxas"
$ws.Cells.Item(3, 2).Value = "DATE"
$ws.Cells.Item(4, 1).Value = "Kamara,"
$ws.Cells.Item(4, 2).Value = "LOCATION"
$ws.Cells.Item(5, 1).Value = "Private Hospitals Association of the Philippines"
$ws.Cells.Item(5, 2).Value = "LOCATION"
$ws.Cells.Item(6, 1).Value = "Arroyo Ad­ministrasyon"
$ws.Cells.Item(6, 2).Value = "PERSON"
$ws.Cells.Item(7, 1).Value = "Arroyo"
$ws.Cells.Item(7, 2).Value = "PERSON"
$ws.Cells.Item(8, 1).Value = "Ayon"
$ws.Cells.Item(8, 2).Value = "PERSON"
$ws.Cells.Item(9, 1).Value = "Charter Change"
$ws.Cells.Item(9, 2).Value = "PERSON"
$ws.Cells.Item(10, 1).Value = "Gayunman,"
$ws.Cells.Item(10, 2).Value = "PERSON"
$ws.Cells.Item(11, 1).Value = "Gloria Forever Constitutionâ?"
$ws.Cells.Item(11, 2).Value = "PERSON"
$ws.Cells.Item(12, 1).Value = "Glo­ria Forever Constitution"
$ws.Cells.Item(12, 2).Value = "PERSON"
$ws.Cells.Item(13, 1).Value = "Ilonggong"
$ws.Cells.Item(13, 2).Value = "PERSON"
$ws.Cells.Item(14, 1).Value = "Ito"
$ws.Cells.Item(14, 2).Value = "PERSON"
$ws.Cells.Item(15, 1).Value = "Konstitus­yon"
$ws.Cells.Item(15, 2).Value = "PERSON"
$ws.Cells.Item(16, 1).Value = "Mas"
$ws.Cells.Item(16, 2).Value = "PERSON"
$ws.Cells.Item(17, 1).Value = "Napunit"
$ws.Cells.Item(17, 2).Value = "PERSON"
$ws.Cells.Item(18, 1).Value = "Pangulong Arroyo"
$ws.Cells.Item(18, 2).Value = "PERSON"
$ws.Cells.Item(19, 1).Value = "Pangu­lo"
$ws.Cells.Item(19, 2).Value = "PERSON"
$ws.Cells.Item(20, 1).Value = "Pilipino"
$ws.Cells.Item(20, 2).Value = "PERSON"
$ws.Cells.Item(21, 1).Value = "Presidente Gloria Arroyo"
$ws.Cells.Item(21, 2).Value = "PERSON"
$ws.Cells.Item(22, 1).Value = "Punit-punit"
$ws.Cells.Item(22, 2).Value = "PERSON"
$ws.Cells.Item(23, 1).Value = "Roxas"
$ws.Cells.Item(23, 2).Value = "PERSON"
$ws.Cells.Item(24, 1).Value = "Roxas,"
$ws.Cells.Item(24, 2).Value = "PERSON"
$ws.Cells.Item(25, 1).Value = "Walang"
$ws.Cells.Item(25, 2).Value = "PERSON"

$ws.Range("A1:B25").Select() | Out-Null
